$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.27000000000035
$ws.Range("H2").Value = [double]"1.319730192719354e-16"
$ws.Range("K2").Value = 47.30517063274595
$ws.Range("L2").Value = "[44.9505454289704, 49.659795836521496]"
$ws.Range("O2").Value = 1.566079220708426
$ws.Range("P2").Value = "[1.515763422452733, 1.616395018964119]"
$ws.Range("S2").Value = 50.90800967283721
$ws.Range("T2").Value = "[49.26892568966874, 52.54709365600568]"
$ws.Range("W2").Value = 18.22072072072098
$ws.Range("X2").Value = 18.02636636636662
$ws.Range("Y2").Value = 18.41507507507534

# Row 3 updates
$ws.Range("E3").Value = 24.94000000000046
$ws.Range("H3").Value = [double]"1.319730192719354e-16"
$ws.Range("K3").Value = 45.54006655327249
$ws.Range("L3").Value = "[42.12998767743372, 48.95014542911126]"
$ws.Range("O3").Value = 2.849132076228581
$ws.Range("P3").Value = "[2.773658378845042, 2.9246057736121203]"
$ws.Range("S3").Value = 49.77429576724634
$ws.Range("T3").Value = "[47.82083511173683, 51.727756422755846]"
$ws.Range("W3").Value = 13.63087087087112
$ws.Range("X3").Value = 13.33129129129153
$ws.Range("Y3").Value = 13.93045045045071
